# SectorGroup.xlsx: the codeforiati:category-code and codeforiati:category-name
# columns (F and G) have had their values swapped for every row (header included),
# i.e. column F now holds the category *name* and column G now holds the
# category *code*.
#
# We swap the two columns via Copy/Paste (through a scratch column) rather than
# reading/writing .Value directly, because re-entering a numeric-looking string
# (e.g. "111") through .Value gets auto-coerced to a number — Copy/Paste
# preserves the original text cell type instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$rows = $ur.Rows.Count()

$colF = $ws.Range($ws.Cells.Item(1, 6), $ws.Cells.Item($rows, 6))
$colG = $ws.Range($ws.Cells.Item(1, 7), $ws.Cells.Item($rows, 7))
$scratch = $ws.Range($ws.Cells.Item(1, 26), $ws.Cells.Item($rows, 26))

$colF.Copy($scratch)
$colG.Copy($colF)
$scratch.Copy($colG)
$scratch.Clear()
